$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,("ECs", "Angpt2", "Tie1", "ECs", 3, 1, 24.62964533333333, 73.888936, 0.7460691037955391, 0.746069103795539, 3, 1, 66.04433166666666, 198.132995, 0.9905709103863496, 0.9905709103863498, 1626.648465227036, 14639.83618704332, 0.7390343513578752, 0.7390343513578752)
    ,("ECs", "Angpt2", "Tie1", "FAPs", 3, 1, 24.62964533333333, 73.888936, 0.7460691037955391, 0.746069103795539, 1, 0.3333333333333333, 0.09667466666666667, 0.290024, 0.001449982309679873, 0.001449982309679873, 2.381062752718222, 21.429564774464, 0.001081787002302249, 0.001081787002302249)
    ,("ECs", "Angpt2", "Tie1", "M2", 3, 1, 24.62964533333333, 73.888936, 0.7460691037955391, 0.746069103795539, 1, 0.3333333333333333, 0.02943566666666667, 0.088307, 0.0004414930758175203, 0.0004414930758175204, 0.7249900301502221, 6.524910271352, 0.0003293843434071133, 0.0003293843434071135)
    ,("ECs", "Angpt2", "Tie1", "sCs", 3, 1, 24.62964533333333, 73.888936, 0.7460691037955391, 0.746069103795539, 3, 1, 0.5025553333333334, 1.507666, 0.007537614228152894, 0.007537614228152896, 12.37775962037511, 111.399836583376, 0.005623581091954534, 0.005623581091954534)
    ,("FAPs", "Angpt2", "Tie1", "ECs", 3, 1, 1.769696333333333, 5.309089, 0.05360677100832464, 0.05360677100832464, 3, 1, 66.04433166666666, 198.132995, 0.9905709103863496, 0.9905709103863498, 116.8784115879506, 1051.905704291555, 0.05310130796058871, 0.05310130796058872)
    ,("FAPs", "Angpt2", "Tie1", "FAPs", 3, 1, 1.769696333333333, 5.309089, 0.05360677100832464, 0.05360677100832464, 1, 0.3333333333333333, 0.09667466666666667, 0.290024, 0.001449982309679873, 0.001449982309679873, 0.1710848031262222, 1.539763228136, 0.0000777288696411306, 0.0000777288696411306)
    ,("FAPs", "Angpt2", "Tie1", "M2", 3, 1, 1.769696333333333, 5.309089, 0.05360677100832464, 0.05360677100832464, 1, 0.3333333333333333, 0.02943566666666667, 0.088307, 0.0004414930758175203, 0.0004414930758175204, 0.05209219136922222, 0.468829722323, 0.00002366701821711072, 0.00002366701821711072)
    ,("FAPs", "Angpt2", "Tie1", "sCs", 3, 1, 1.769696333333333, 5.309089, 0.05360677100832464, 0.05360677100832464, 3, 1, 0.5025553333333334, 1.507666, 0.007537614228152894, 0.007537614228152896, 0.8893703306971111, 8.004332976274, 0.0004040671598776819, 0.0004040671598776819)
    ,("M2", "Angpt2", "Tie1", "ECs", 3, 1, 1.275024, 3.825072, 0.03862239996247084, 0.03862239996247083, 3, 1, 66.04433166666666, 198.132995, 0.9905709103863496, 0.9905709103863498, 84.20810793896, 757.87297145064, 0.03825822589213045, 0.03825822589213046)
    ,("M2", "Angpt2", "Tie1", "FAPs", 3, 1, 1.275024, 3.825072, 0.03862239996247084, 0.03862239996247083, 1, 0.3333333333333333, 0.09667466666666667, 0.290024, 0.001449982309679873, 0.001449982309679873, 0.123262520192, 1.109362681728, 0.0000560017967029633, 0.0000560017967029633)
    ,("M2", "Angpt2", "Tie1", "M2", 3, 1, 1.275024, 3.825072, 0.03862239996247084, 0.03862239996247083, 1, 0.3333333333333333, 0.02943566666666667, 0.088307, 0.0004414930758175203, 0.0004414930758175204, 0.037531181456, 0.337780633104, 0.00001705152215488573, 0.00001705152215488574)
    ,("M2", "Angpt2", "Tie1", "sCs", 3, 1, 1.275024, 3.825072, 0.03862239996247084, 0.03862239996247083, 3, 1, 0.5025553333333334, 1.507666, 0.007537614228152894, 0.007537614228152896, 0.640770111328, 5.766931001952, 0.000291120751482532, 0.000291120751482532)
    ,("sCs", "Angpt2", "Tie1", "ECs", 3, 1, 5.338186666666666, 16.01456, 0.1617017252336654, 0.1617017252336654, 3, 1, 66.04433166666666, 198.132995, 0.9905709103863496, 0.9905709103863498, 352.556970711911, 3173.0127364072, 0.1601770251757553, 0.1601770251757553)
    ,("sCs", "Angpt2", "Tie1", "FAPs", 3, 1, 5.338186666666666, 16.01456, 0.1617017252336654, 0.1617017252336654, 1, 0.3333333333333333, 0.09667466666666667, 0.290024, 0.001449982309679873, 0.001449982309679873, 0.5160674166044444, 4.64460674944, 0.0002344646410335303, 0.0002344646410335303)
    ,("sCs", "Angpt2", "Tie1", "M2", 3, 1, 5.338186666666666, 16.01456, 0.1617017252336654, 0.1617017252336654, 1, 0.3333333333333333, 0.02943566666666667, 0.088307, 0.0004414930758175203, 0.0004414930758175204, 0.1571330833244444, 1.41419774992, 0.00007139019203841048, 0.00007139019203841049)
    ,("sCs", "Angpt2", "Tie1", "sCs", 3, 1, 5.338186666666666, 16.01456, 0.1617017252336654, 0.1617017252336654, 3, 1, 0.5025553333333334, 1.507666, 0.007537614228152894, 0.007537614228152896, 2.682734179662222, 24.14460761696, 0.001218845224838146, 0.001218845224838146)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = 2 + $i
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}
